# Update the Excel workbook "Prix Spot" sheet with a new day column (R)
# containing the 01-jul prices, mirroring the existing 15-jun..30-jun
# columns (B..Q).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last existing header cell (Q1) onto the new
# header cell (R1) so it reuses the same header style (bold, bordered,
# centered) instead of creating a brand new style entry.
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R1").Value = "01-jul"

# New numeric values for the 01-jul column (rows 2-25).
$values = @(
    111.28,
    95.41,
    89.09999999999999,
    88.08,
    90.01000000000001,
    96.64,
    114.97,
    120.06,
    114.74,
    93.19,
    94.91,
    80.7,
    75.11,
    65.53,
    68.95,
    78.16,
    90.8,
    102.65,
    150,
    223.54,
    235,
    215,
    199.3,
    129.78
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 18).Value = $values[$i]
}
